# Insert a new row at position 9 (shifts existing rows 9..110 down to 10..111)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with its data
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 45163
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112035
$ws.Range("G9").Value = "Bruselas (repollito)"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 350
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 18686
$ws.Range("N9").Value = "$/malla 15 kilos"
$ws.Range("O9").Value = "Provincia de Quillota"
$ws.Range("P9").Value = 1246
$ws.Range("Q9").Value = 15
$ws.Range("R9").Value = "Hortaliza"
